$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sequence Diagram")
$ws.Range("B25").Value = "test"
